# Update attendance/interest numbers ("想去人数") for several events.
# Sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5359
$ws1.Range("F3").Value = 583
$ws1.Range("F4").Value = 11266
$ws1.Range("G4").Value = 58
$ws1.Range("F6").Value = 585
$ws1.Range("F7").Value = 160
$ws1.Range("F8").Value = 241
$ws1.Range("F9").Value = 969
$ws1.Range("F10").Value = 94

# Sheet "全部类型" (All types) mirrors the same events at different rows.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5359
$ws4.Range("F5").Value = 583
$ws4.Range("F7").Value = 11266
$ws4.Range("G7").Value = 58
$ws4.Range("F9").Value = 585
$ws4.Range("F10").Value = 160
$ws4.Range("F13").Value = 241
$ws4.Range("F14").Value = 969
$ws4.Range("F16").Value = 94
